{"js": "// The \"Generate report\" list of childless tags was being appended twice\n// to the document (a duplicate bug). This removes the first (duplicate)\n// copy of the list, keeping only one copy, right after the intro line\n// \"These are the childless tags that were found in the documents: \".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst texts = paragraphs.items.map((p) => p.text);\nconst total = texts.length;\n\n// Locate two adjacent, equal-length runs of consecutive paragraphs that\n// are identical (the duplicated tag list) appearing after the heading\n// paragraphs. We search for the longest such duplicated run.\nlet blockLength = 0;\nlet firstStart = -1;\nlet secondStart = -1;\n\nsearchLengths:\nfor (let len = Math.floor((total - 2) / 2); len >= 1; len--) {\n  for (let i = 2; i + 2 * len <= total; i++) {\n    const j = i + len;\n    let isMatch = true;\n    for (let k = 0; k < len; k++) {\n      if (texts[i + k] !== texts[j + k]) {\n        isMatch = false;\n        break;\n      }\n    }\n    if (isMatch) {\n      blockLength = len;\n      firstStart = i;\n      secondStart = j;\n      break searchLengths;\n    }\n  }\n}\n\nif (blockLength > 0) {\n  // Delete the first occurrence of the duplicated block (the earlier\n  // copy), leaving the second copy intact. Delete from the end of the\n  // range backward so indices of earlier items remain valid.\n  for (let idx = firstStart + blockLength - 1; idx >= firstStart; idx--) {\n    paragraphs.items[idx].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# The \"Generate report\" list of childless tags was being appended twice\n# to the document (a duplicate bug caused by the Generate report button\n# not being disabled while it ran). This removes the first (duplicate)\n# copy of the list, keeping only one copy, right after the intro line\n# \"These are the childless tags that were found in the documents: \".\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n# Collect paragraph text (without the trailing paragraph-mark character).\n$texts = @()\nfor ($i = 1; $i -le $count; $i++) {\n    $raw = $d.Paragraphs.Item($i).Range.Text\n    $texts += $raw.TrimEnd([char]13)\n}\n\n$n = $texts.Length\n\n# Locate two adjacent, equal-length runs of consecutive paragraphs that are\n# identical (the duplicated tag list), searching for the longest such run\n# starting after the first two heading paragraphs.\n$blockLen = 0\n$firstStart = -1\n$secondStart = -1\n\n$maxLen = [Math]::Floor(($n - 2) / 2)\nfor ($len = $maxLen; $len -ge 1; $len--) {\n    $found = $false\n    for ($i = 2; ($i + 2 * $len) -le $n; $i++) {\n        $j = $i + $len\n        $match = $true\n        for ($k = 0; $k -lt $len; $k++) {\n            if ($texts[$i + $k] -ne $texts[$j + $k]) { $match = $false; break }\n        }\n        if ($match) {\n            $blockLen = $len\n            $firstStart = $i\n            $secondStart = $j\n            $found = $true\n            break\n        }\n    }\n    if ($found) { break }\n}\n\nif ($blockLen -gt 0) {\n    # Convert the 0-based paragraph index of the first duplicated run to\n    # Word's 1-based Paragraphs collection indexing, then delete the whole\n    # run (from the start of its first paragraph to the end of its last\n    # paragraph) in a single range delete.\n    $firstParaIndex = $firstStart + 1\n    $lastParaIndex = $firstStart + $blockLen\n\n    $startPos = $d.Paragraphs.Item($firstParaIndex).Range.Start\n    $endPos = $d.Paragraphs.Item($lastParaIndex).Range.End\n\n    $r = $d.Range($startPos, $endPos)\n    $r.Delete()\n}\n"}
